$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, shifting existing rows 149-229 down to 150-230.
$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with the new record's data.
$ws.Range("A149").Value = 8
$ws.Range("B149").Value = "Terminal La Palmera de La Serena"
$ws.Range("C149").Value = "Coquimbo"
$ws.Range("D149").Value = 44455
$ws.Range("E149").Value = 4
$ws.Range("F149").Value = 100114001
$ws.Range("G149").Value = "Papa"
$ws.Range("H149").Value = "Cardinal"
$ws.Range("I149").Value = "1a (cosecha)"
$ws.Range("J149").Value = 2000
$ws.Range("K149").Value = 11500
$ws.Range("L149").Value = 12000
$ws.Range("M149").Value = 11750
$ws.Range("N149").Value = '$/saco 25 kilos'
$ws.Range("O149").Value = "Provincia del Elquí"
$ws.Range("P149").Value = 470
$ws.Range("Q149").Value = 25
$ws.Range("R149").Value = "Hortaliza"
